$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ShipmentTracking (P), ActualRate (Q) and Result (R) columns hold
# numeric-looking values ("320018256721", "$223.37", ...) that must be
# stored as plain text (shared strings), matching every other cell already
# in these columns. Forcing NumberFormat to "@" before the assignment
# prevents Excel from auto-converting the values to numbers, and resetting
# the style back to "Normal" afterwards keeps the cells unstyled (style
# index 0), same as the surrounding, untouched cells.
$pRange = $ws.Range("P2:P25")
$qRange = $ws.Range("Q22:Q25")
$rRange = $ws.Range("R22:R25")

$pRange.NumberFormat = "@"
$qRange.NumberFormat = "@"
$rRange.NumberFormat = "@"

# New FedEx tracking numbers (ShipmentTracking) for the already-quoted rows.
$ws.Range("P2").Value = '320018256721'
$ws.Range("P3").Value = '320018256732'
$ws.Range("P4").Value = '320018256765'
$ws.Range("P5").Value = '320018256787'
$ws.Range("P6").Value = '320018256824'
$ws.Range("P7").Value = '320018256846'
$ws.Range("P8").Value = '320018256879'
$ws.Range("P9").Value = '320018256890'
$ws.Range("P10").Value = '320018256927'
$ws.Range("P11").Value = '320018256949'
$ws.Range("P12").Value = '320018256982'
$ws.Range("P13").Value = '320018257007'
$ws.Range("P14").Value = '320018257030'
$ws.Range("P15").Value = '320018257051'
$ws.Range("P16").Value = '320018257084'
$ws.Range("P17").Value = '320018257100'
$ws.Range("P18").Value = '320018257143'
$ws.Range("P19").Value = '320018257165'
$ws.Range("P20").Value = '320018257198'
$ws.Range("P21").Value = '320018257213'

# Rows 22-25 got an updated quote: new tracking number, new actual rate and
# the result flips from PASS to FAIL.
$ws.Range("P22").Value = '320018257246'
$ws.Range("Q22").Value = '$223.37'
$ws.Range("R22").Value = 'FAIL'

$ws.Range("P23").Value = '320018257257'
$ws.Range("Q23").Value = '$436.98'
$ws.Range("R23").Value = 'FAIL'

$ws.Range("P24").Value = '320018257268'
$ws.Range("Q24").Value = '$278.12'
$ws.Range("R24").Value = 'FAIL'

$ws.Range("P25").Value = '320018257279'
$ws.Range("Q25").Value = '$52.88'
$ws.Range("R25").Value = 'FAIL'

$pRange.Style = "Normal"
$qRange.Style = "Normal"
$rRange.Style = "Normal"
